# Removed Extension Payments Tax Type from execution.
#
# The test run that previously exercised the "Extension Payments" row (row 3)
# is excluded by flipping its Execute flag from "Y" to "DONOTRUN". The
# timestamps for the rows that DID run (rows 2 and 4) are bumped to the
# latest execution pass, and row 3 keeps an older timestamp since it no
# longer executes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Thu Dec 07 21:44:04 EST 2023"
$ws.Range("B3").Value = "Wed Nov 01 15:54:52 EDT 2023"
$ws.Range("C3").Value = "DONOTRUN"
$ws.Range("B4").Value = "Thu Dec 07 21:44:17 EST 2023"

# Column C was widened (and is no longer auto "best fit") to comfortably
# show "DONOTRUN".
$ws.Columns("C").ColumnWidth = 13.8

# The cell that was actually edited (C3) is left as the active selection.
$ws.Range("C3").Select()
